$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.695.73"
$ws.Range("E2").Value = "  -2.98%  "
$ws.Range("D3").Value = "2.095.41"
$ws.Range("E3").Value = "  -2.24%  "
$ws.Range("D4").Value = "'1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'345.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.05%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "'0.5156"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.97%  "
$ws.Range("E8").Value = "  -3.76%  "
$ws.Range("D9").Value = "'52.53"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.09237"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("E11").Value = "  -1.22%  "
$ws.Range("D12").Value = "'24.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.00%  "
$ws.Range("D13").Value = "2.088.61"
$ws.Range("E13").Value = "  -2.33%  "
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("D16").Value = "'99.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.69%  "
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "'20.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.46%  "
$ws.Range("D20").Value = "'0.06658"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("E22").Value = "  -2.55%  "
$ws.Range("D23").Value = "29.750.28"
$ws.Range("E23").Value = "  -3.11%  "
$ws.Range("D24").Value = "'12.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.73%  "
$ws.Range("D25").Value = "'2.319"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.58%  "
$ws.Range("D26").Value = "2.337.65"
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("E27").Value = "  -2.88%  "
$ws.Range("D28").Value = "'2.525"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.05%  "
$ws.Range("D29").Value = "'162.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").Value = "'133.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.65%  "
$ws.Range("E31").Value = "  -7.39%  "
$ws.Range("E32").Value = "  -2.76%  "
$ws.Range("E33").Value = "  -1.44%  "
$ws.Range("D34").Value = "'6.178"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.24%  "
$ws.Range("D35").Value = "'3.934"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.08%  "
$ws.Range("D36").Value = "'6.202"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").Value = "'10.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.89%  "
$ws.Range("D38").Value = "'0.02569"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.28%  "
$ws.Range("D39").Value = "'0.06701"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.76%  "
$ws.Range("E40").Value = "  -2.28%  "
$ws.Range("D41").Value = "'0.6862"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.48%  "
$ws.Range("D42").Value = "'0.2227"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.96%  "
$ws.Range("E43").Value = "  +1.59%  "
$ws.Range("D44").Value = "'0.6641"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.45%  "
$ws.Range("D45").Value = "'14.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.81%  "
$ws.Range("D46").Value = "'2.314"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.27%  "
$ws.Range("D47").Value = "'3.626"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.50%  "
$ws.Range("D48").Value = "'0.00000000352"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.69%  "
$ws.Range("E49").Value = "  -3.03%  "
$ws.Range("D50").Value = "'82.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("D51").Value = "'0.3307"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.58%  "
